$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 578
$ws1.Range("F11").Value = 81
$ws1.Range("F12").Value = 11438
$ws1.Range("F13").Value = 6504
$ws1.Range("F20").Value = 904
$ws1.Range("F22").Value = 254
$ws1.Range("F35").Value = 219

# Sheet "演出" (Performance) - sheet2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 3659

# Sheet "全部类型" (All types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 2748
$ws4.Range("F13").Value = 578
$ws4.Range("F18").Value = 81
$ws4.Range("F19").Value = 11438
$ws4.Range("F20").Value = 3659
$ws4.Range("F21").Value = 6504
$ws4.Range("F28").Value = 904
$ws4.Range("F43").Value = 219
